$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as text, matching the source data
# (many of these values look numeric, e.g. "23.10" or "0.9999", and would otherwise
# be auto-converted to numbers by Excel, stripping formatting / trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.164.32'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '1.827.56'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '241.71'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').Value = '0.6237'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.07333'
$ws.Range('E8').Value = '  -1.11%  '
$ws.Range('D9').Value = '0.2906'
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('D10').Value = '23.10'
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('D11').Value = '0.07679'
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').Value = '1.837.77'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '4.953'
$ws.Range('E13').Value = '  -1.10%  '
$ws.Range('D14').Value = '0.6655'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').Value = '82.27'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').Value = '0.000008896'
$ws.Range('E16').Value = '  -4.17%  '
$ws.Range('D17').Value = '5.846'
$ws.Range('E17').Value = '  -1.63%  '
$ws.Range('D18').Value = '29.142.98'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = '2.076.57'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').Value = '238.14'
$ws.Range('E20').Value = '  +6.96%  '
$ws.Range('D21').Value = '12.44'
$ws.Range('E21').Value = '  -1.35%  '
$ws.Range('D22').Value = '0.9995'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').Value = '7.333'
$ws.Range('E23').Value = '  +3.23%  '
$ws.Range('D24').Value = '0.9999'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').Value = '158.04'
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('D26').Value = '0.1424'
$ws.Range('E26').Value = '  +2.16%  '
$ws.Range('D27').Value = '8.494'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').Value = '17.67'
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('D29').Value = '1.482'
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('D30').Value = '0.05749'
$ws.Range('E30').Value = '  +1.06%  '
$ws.Range('D31').Value = '4.091'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').Value = '4.095'
$ws.Range('E32').Value = '  -1.56%  '
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('D34').Value = '1.865'
$ws.Range('E34').Value = '  +2.41%  '
$ws.Range('D35').Value = '0.7315'
$ws.Range('E35').Value = '  -1.26%  '
$ws.Range('D36').Value = '1.134'
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('D37').Value = '2.625'
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('D38').Value = '2.838'
$ws.Range('E38').Value = '  +2.79%  '
$ws.Range('D39').Value = '1.218.36'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = '6.315'
$ws.Range('E41').Value = '  -2.94%  '
$ws.Range('D42').Value = '0.9218'
$ws.Range('E42').Value = '  +3.77%  '
$ws.Range('D43').Value = '0.9999'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '101.85'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').Value = '1.984.55'
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').Value = '64.77'
$ws.Range('E46').Value = '  -2.35%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.00000000121'
$ws.Range('E47').Value = '  -2.61%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.5085'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('D49').Value = '0.4035'
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('D50').Value = '9.133'
$ws.Range('E50').Value = '  +0.83%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.1129'
$ws.Range('E51').Value = '  +2.65%  '
